$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header counts (sample sizes) changed
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON) values tweaked; C2 deleted
$ws.Range("B2").Value = 9.0122213228944847
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = -0.12996849305077129
$ws.Range("E2").Value = -0.67359004972923542

# Row 3 (STR) values tweaked; B3 deleted
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 3.3215006037958563
$ws.Range("D3").Value = 2.1252554919599076
$ws.Range("E3").Value = -2.3892836560151847

# Selection now only covers the updated data block
$ws.Range("B1:E3").Select()
